# Refresh cryptos list with updated prices / volume figures (scraped Sat Oct 26 2024).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "66.887.47"
$ws.Range('E2').Value = "  -1.91%  "
$ws.Range('D3').Value = "2.451.52"
$ws.Range('E3').Value = "  -3.12%  "
$ws.Range('D4').Value = "'0.999"
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = "  -0.07%  "
$ws.Range('D5').Value = "'578.54"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = "  -3.05%  "
$ws.Range('D6').Value = "'165.16"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = "  -6.12%  "
$ws.Range('E7').Value = "  +0.05%  "
$ws.Range('D8').Value = "'0.511"
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = "  -3.64%  "
$ws.Range('D9').Value = "2.452.38"
$ws.Range('E9').Value = "  -3.02%  "
$ws.Range('D10').Value = "'0.134"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = "  -4.87%  "
$ws.Range('D11').Value = "'0.164"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = "  -1.02%  "
$ws.Range('B12').Value = "Cardano"
$ws.Range('C12').Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range('D12').Value = "'0.331"
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = "  -4.05%  "
$ws.Range('B13').Value = "Toncoin"
$ws.Range('C13').Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range('D13').Value = "'4.86"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = "  -4.98%  "
$ws.Range('B14').Value = "Avalanche"
$ws.Range('C14').Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range('D14').Value = "'25.29"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = "  -5.26%  "
$ws.Range('B15').Value = "WrappedliquidstakedEther2.0"
$ws.Range('C15').Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range('D15').Value = "2.888.20"
$ws.Range('E15').Value = "  -3.42%  "
$ws.Range('D16').Value = "66.629.11"
$ws.Range('E16').Value = "  -1.85%  "
$ws.Range('D17').Value = "'0.0000168"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = "  -6.38%  "
$ws.Range('D18').Value = "2.470.64"
$ws.Range('E18').Value = "  -0.89%  "
$ws.Range('D19').Value = "'11.37"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = "  -5.41%  "
$ws.Range('D20').Value = "'7.76"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = "  -4.80%  "
$ws.Range('D21').Value = "'353.69"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = "  -3.21%  "
$ws.Range('D22').Value = "'4.05"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = "  -3.03%  "
$ws.Range('E23').Value = "  +0.13%  "
$ws.Range('D24').Value = "'69.33"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = "  -2.80%  "
$ws.Range('D25').Value = "'4.20"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = "  -10.06%  "
$ws.Range('D26').Value = "'1.75"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = "  -9.14%  "
$ws.Range('D27').Value = "'8.91"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = "  -11.67%  "
$ws.Range('D28').Value = "'1.00"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = "  +0.05%  "
$ws.Range('D29').Value = "2.568.69"
$ws.Range('E29').Value = "  -2.83%  "
$ws.Range('D30').Value = "0.0₃0898"
$ws.Range('E30').Value = "  -9.14%  "
$ws.Range('D31').Value = "'504.57"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = "  -5.64%  "
$ws.Range('D32').Value = "'7.79"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = "  -6.96%  "
$ws.Range('D33').Value = "'1.77"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = "  -7.34%  "
$ws.Range('D34').Value = "'1.22"
$ws.Range('D34').ClearFormats()
$ws.Range('D35').Value = "'0.999"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = "  -0.10%  "
$ws.Range('D36').Value = "'159.02"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = "  +0.89%  "
$ws.Range('D37').Value = "'0.117"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = "  -9.86%  "
$ws.Range('D38').Value = "'18.56"
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = "  -0.70%  "
$ws.Range('D39').Value = "'18.44"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = "  -1.92%  "
$ws.Range('D40').Value = "'1.34"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = "  -8.07%  "
$ws.Range('E41').Value = "  -0.03%  "
$ws.Range('E42').Value = "  -7.75%  "
$ws.Range('D43').Value = "'0.326"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = "  -7.40%  "
$ws.Range('D44').Value = "'4.72"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = "  -9.26%  "
$ws.Range('E45').Value = "  -3.29%  "
$ws.Range('E46').Value = "  -9.95%  "
$ws.Range('D47').Value = "'141.42"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = "  -4.48%  "
$ws.Range('D48').Value = "'3.47"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = "  -6.80%  "
$ws.Range('D49').Value = "'0.513"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = "  -8.05%  "
$ws.Range('D50').Value = "'1.59"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = "  -8.70%  "
$ws.Range('D51').Value = "'0.0731"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = "  -2.84%  "
